# Update the Test Case Id column (A) with the new "UT BTS TC_*" naming scheme.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value  = "UT- BTS TC_1"
$ws.Range("A8").Value  = "UT BTS TC_1.2"
$ws.Range("A9").Value  = "UT BTS TC_2.1"
$ws.Range("A10").Value = "UT BTS TC_2.3"
$ws.Range("A11").Value = "UT BTS TC_2.2"
$ws.Range("A13").Value = "UT BTS TC_3.1"
$ws.Range("A14").Value = "UT BTS TC_3.2"
$ws.Range("A15").Value = "UT BTS TC_3.3"
$ws.Range("A16").Value = "UT BTS TC_3.4"
$ws.Range("A17").Value = "UT BTS TC_4"

# Move the active selection to A17, matching the saved cursor position.
$ws.Range("A17").Select()
